# This script reproduces a weekly data refresh: a new daily price record
# (dated 44827) is inserted as a new row right before the existing row for
# date 44721, pushing all subsequent rows (old rows 69-138) down by one.
# The workbook's dimension grows from A1:R138 to A1:R139 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 69; Excel automatically shifts the
# existing rows 69..138 down to 70..139 and extends the used range.
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new record's data. All
# "constant" columns (mercado, region, categoria, variedad, calidad,
# unidad, origen, kg/unidades, clasificacion) match the surrounding rows
# for this product/market, while the date, volume and price columns carry
# the new values.
$ws.Range("A69").Value2 = 10
$ws.Range("B69").Value2 = "Vega Modelo de Temuco"
$ws.Range("C69").Value2 = "La Araucanía"
$ws.Range("D69").Value2 = 44827
$ws.Range("E69").Value2 = 9
$ws.Range("F69").Value2 = 100112035
$ws.Range("G69").Value2 = "Bruselas (repollito)"
$ws.Range("H69").Value2 = "Sin especificar"
$ws.Range("I69").Value2 = "Primera"
$ws.Range("J69").Value2 = 30
$ws.Range("K69").Value2 = 24000
$ws.Range("L69").Value2 = 24000
$ws.Range("M69").Value2 = 24000
$ws.Range("N69").Value2 = "`$/malla 10 kilos"
$ws.Range("O69").Value2 = "Región Metropolitana"
$ws.Range("P69").Value2 = 2400
$ws.Range("Q69").Value2 = 10
$ws.Range("R69").Value2 = "Hortaliza"
